$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "ShearF"

# 2. Tiny floating point corrections in row 13 (last-bit precision refresh)
$ws.Range("D13").Value = 0.9933785644704798
$ws.Range("J13").Value = 0.9933785644704798
$ws.Range("K13").Value = 0.9933077201834214
$ws.Range("L13").Value = 0.9950650458092126

# 3. Tiny floating point corrections in row 15
$ws.Range("D15").Value = 0.729469432657212
$ws.Range("J15").Value = 0.729469432657212
$ws.Range("K15").Value = 0.8819903124027868

# 4. Append new row 16 with results for HKL index 14 / HexGrid-60degTilt5degRes
#    Copy formatting from row 15's A cell (bordered/bold/centered style) first
$ws.Range("A15").Copy($ws.Range("A16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.225040886154097
$ws.Range("D16").Value = 2.34493493238325
$ws.Range("E16").Value = 0.6296575682396935
$ws.Range("F16").Value = 1.225040886154097
$ws.Range("G16").Value = 1.332156129456757
$ws.Range("H16").Value = 0.4592399385408076
$ws.Range("I16").Value = 0.7660167035182626
$ws.Range("J16").Value = 2.34493493238325
$ws.Range("K16").Value = 1.487296250311472
$ws.Range("L16").Value = 1.356168568232785
$ws.Range("M16").Value = 1.126174359715478
